# Leave Card update — 2018 LEAVE CREDITS & 2017 LEAVE BALANCE sheets
# (commit: "Update Leave Card 12/22/2023 10:59 AM")

$wb = $excel.ActiveWorkbook

$ws2018 = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws2017 = $wb.Worksheets.Item("2017 LEAVE BALANCE")

# ---------------------------------------------------------------------------
# 2018 LEAVE CREDITS: shift the monthly "PERIOD" dates from the 1st of each
# month to the last day of each month (rows 77-95), and post the 1.25-day
# VL/SL credit that had been missing for Apr-Oct 2023 (rows 80-86).
# ---------------------------------------------------------------------------

$ws2018.Range("A77").Value2 = 44957
$ws2018.Range("A78").Value2 = 44985
$ws2018.Range("A79").Value2 = 45016

$ws2018.Range("A80").Value2 = 45046
$ws2018.Range("C80").Value2 = 1.25

$ws2018.Range("A81").Value2 = 45077
$ws2018.Range("C81").Value2 = 1.25

$ws2018.Range("A82").Value2 = 45107
$ws2018.Range("C82").Value2 = 1.25

$ws2018.Range("A83").Value2 = 45138
$ws2018.Range("C83").Value2 = 1.25

$ws2018.Range("A84").Value2 = 45169
$ws2018.Range("C84").Value2 = 1.25

$ws2018.Range("A85").Value2 = 45199
$ws2018.Range("C85").Value2 = 1.25

$ws2018.Range("A86").Value2 = 45230
$ws2018.Range("C86").Value2 = 1.25

$ws2018.Range("A87").Value2 = 45260
$ws2018.Range("A88").Value2 = 45291
$ws2018.Range("A89").Value2 = 45322
$ws2018.Range("A90").Value2 = 45351
$ws2018.Range("A91").Value2 = 45382
$ws2018.Range("A92").Value2 = 45412
$ws2018.Range("A93").Value2 = 45443
$ws2018.Range("A94").Value2 = 45473
$ws2018.Range("A95").Value2 = 45504

# ---------------------------------------------------------------------------
# 2017 LEAVE BALANCE: post a new VL(6-0-0) leave entry on row 51
# ---------------------------------------------------------------------------

$ws2017.Range("A51").Value2 = 45231
$ws2017.Range("B51").Value2 = "VL(6-0-0)"
$ws2017.Range("D51").Value2 = 6
$ws2017.Range("K51").Value2 = "11/22-24,27-29/2023"

# ---------------------------------------------------------------------------
# Window / view state: the workbook was left with the "2017 LEAVE BALANCE"
# tab active, scrolled so row 40 (resp. row 73 on 2018 LEAVE CREDITS) is the
# first row under the split, with I9 selected in the header pane and the
# newly-edited cell selected in the lower pane.
# ---------------------------------------------------------------------------

$ws2018.Activate()
$excel.ActiveWindow.SplitRow = 72
$ws2018.Range("I9").Select()
$ws2018.Range("G89").Select()

$ws2017.Activate()
$excel.ActiveWindow.SplitRow = 39
$ws2017.Range("I9").Select()
$ws2017.Range("K51").Select()
